$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.407.70"
$ws.Range("E2").Value = "  +4.41%  "
$ws.Range("D3").Value = "1.727.45"
$ws.Range("E3").Value = "  +2.38%  "
$ws.Range("E4").Value = "  -0.18%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "218.75"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("E6").Value = "  +0.44%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"
$ws.Range("E7").Value = "  -0.24%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "24.00"
$r.Style = "Normal"
$ws.Range("E8").Value = "  +3.28%  "
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("E10").Value = "  +1.12%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0893"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "1.971.64"
$ws.Range("D13").Value = "1.726.61"
$ws.Range("E13").Value = "  +2.31%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "4.24"
$r.Style = "Normal"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("E15").Value = "  +1.73%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "67.58"
$r.Style = "Normal"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "28.365.44"
$ws.Range("E17").Value = "  +4.26%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "246.66"
$r.Style = "Normal"
$ws.Range("E18").Value = "  +4.40%  "
$ws.Range("D19").Value = "0.0₃0751"
$ws.Range("E19").Value = "  +0.78%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "7.90"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -2.09%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  +0.89%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "9.65"
$r.Style = "Normal"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  -1.42%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "149.49"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("E26").Value = "  +1.79%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "16.61"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("E29").Value = "  -0.30%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "0.0516"
$r.Style = "Normal"
$ws.Range("E30").Value = "  +2.52%  "
$ws.Range("E31").Value = "  +2.94%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").Value = "1.487.93"
$ws.Range("E33").Value = "  -4.14%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  -1.53%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "0.979"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +3.16%  "
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("E39").Value = "  +1.51%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "1.06"
$r.Style = "Normal"
$ws.Range("E40").Value = "  +0.12%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "69.77"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +0.78%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -0.23%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "5.66"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("D44").Value = "1.876.41"
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("E46").Value = "  +2.09%  "
$ws.Range("E47").Value = "  +6.73%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "90.49"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("E49").Value = "  +2.94%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "8.13"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -2.26%  "
$ws.Range("E51").Value = "  -0.89%  "
